# Update 2p3. Added templates for formula student suspension, torque
# vectoring, four-wheel steering.
#
# The existing "Sedan_HambaLG" BodyGeometry template sheet is duplicated to
# create a new "FSAE_Achilles" template, placed immediately after the
# original tab and made the active sheet. Two values are then tweaked on
# the new sheet: the "Instance" label (H3) and the rWheelCutout z/scalar
# value (H6).

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("Sedan_HambaLG")

# Copy the sheet, placing the new copy right after the source sheet. This
# also makes the new sheet the active tab (matching the workbook's saved
# "activeTab" pointing at the 2nd sheet).
$sourceSheet.Copy($null, $sourceSheet)

$newSheet = $wb.Worksheets.Item($sourceSheet.Index + 1)
$newSheet.Name = "FSAE_Achilles"

# Instance label for the new template.
$newSheet.Range("H3").Value = "FSAE_Achilles"

# rWheelCutout z/scalar value differs on the new template.
$newSheet.Range("H6").Value = 0.25
